$d = $word.ActiveDocument

$d.Content.Find.Execute("827×5=4135", $true, $false, $false, $false, $false, $true, 1, $false, "124×9=1116", 2) | Out-Null
$d.Content.Find.Execute("949×9=8541", $true, $false, $false, $false, $false, $true, 1, $false, "509×7=3563", 2) | Out-Null
$d.Content.Find.Execute("232×8=1856", $true, $false, $false, $false, $false, $true, 1, $false, "160×3=480", 2) | Out-Null
$d.Content.Find.Execute("161×3=483", $true, $false, $false, $false, $false, $true, 1, $false, "330×8=2640", 2) | Out-Null
$d.Content.Find.Execute("861×6=5166", $true, $false, $false, $false, $false, $true, 1, $false, "639×2=1278", 2) | Out-Null
$d.Content.Find.Execute("397×8=3176", $true, $false, $false, $false, $false, $true, 1, $false, "888×9=7992", 2) | Out-Null
$d.Content.Find.Execute("751×6=4506", $true, $false, $false, $false, $false, $true, 1, $false, "426×8=3408", 2) | Out-Null
$d.Content.Find.Execute("344×9=3096", $true, $false, $false, $false, $false, $true, 1, $false, "854×6=5124", 2) | Out-Null
$d.Content.Find.Execute("516×2=1032", $true, $false, $false, $false, $false, $true, 1, $false, "882×9=7938", 2) | Out-Null
$d.Content.Find.Execute("859×5=4295", $true, $false, $false, $false, $false, $true, 1, $false, "598×9=5382", 2) | Out-Null
$d.Content.Find.Execute("149×5=745", $true, $false, $false, $false, $false, $true, 1, $false, "263×4=1052", 2) | Out-Null
$d.Content.Find.Execute("722×3=2166", $true, $false, $false, $false, $false, $true, 1, $false, "648×9=5832", 2) | Out-Null
$d.Content.Find.Execute("391×4=1564", $true, $false, $false, $false, $false, $true, 1, $false, "224×4=896", 2) | Out-Null
$d.Content.Find.Execute("844×2=1688", $true, $false, $false, $false, $false, $true, 1, $false, "549×2=1098", 2) | Out-Null
$d.Content.Find.Execute("302×4=1208", $true, $false, $false, $false, $false, $true, 1, $false, "370×7=2590", 2) | Out-Null
$d.Content.Find.Execute("565×4=2260", $true, $false, $false, $false, $false, $true, 1, $false, "934×6=5604", 2) | Out-Null
$d.Content.Find.Execute("217×3=651", $true, $false, $false, $false, $false, $true, 1, $false, "541×9=4869", 2) | Out-Null
$d.Content.Find.Execute("891×5=4455", $true, $false, $false, $false, $false, $true, 1, $false, "537×8=4296", 2) | Out-Null
$d.Content.Find.Execute("476×9=4284", $true, $false, $false, $false, $false, $true, 1, $false, "673×3=2019", 2) | Out-Null
$d.Content.Find.Execute("270×3=810", $true, $false, $false, $false, $false, $true, 1, $false, "199×7=1393", 2) | Out-Null
$d.Content.Find.Execute("307×6=1842", $true, $false, $false, $false, $false, $true, 1, $false, "716×9=6444", 2) | Out-Null
$d.Content.Find.Execute("810×5=4050", $true, $false, $false, $false, $false, $true, 1, $false, "149×4=596", 2) | Out-Null
$d.Content.Find.Execute("175×9=1575", $true, $false, $false, $false, $false, $true, 1, $false, "383×6=2298", 2) | Out-Null
$d.Content.Find.Execute("581×4=2324", $true, $false, $false, $false, $false, $true, 1, $false, "487×7=3409", 2) | Out-Null
$d.Content.Find.Execute("188×5=940", $true, $false, $false, $false, $false, $true, 1, $false, "514×6=3084", 2) | Out-Null
